$d = $word.ActiveDocument

# 1) Merge "proje" + "c" + "t" -> "project" (inside the Tapit Co. hyperlink)
#    Anchor on the unique "Tapit Co.(" prefix, collapse to just after it, then
#    retype "project" over itself so the run-splitting inside the hyperlink is
#    normalised away without touching the text itself.
$anchor = $d.Content
$anchor.Find.Execute("Tapit Co.(", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor.Collapse(0) | Out-Null
$anchor.Find.Execute("project", $true, $false, $false, $false, $false, $true, 1, $false, "project", 2) | Out-Null

# 2) "Telegram, Twitch, and Discord" -> "Twitch and Discord"
$d.Content.Find.Execute("Telegram, Twitch, and Discord", $true, $false, $false, $false, $false, $true, 1, $false, "Twitch and Discord", 2) | Out-Null
